$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new registrant replaces the old bulk-register placeholder row ---
$ws.Range("A2").Value = "bocaioandoru12+5@gmail.com"
$ws.Range("B2").Value = "Doru5"

# --- Move "Department"/"Title" values from J/K into H/I (Campus/Room dropped) ---
$deptHeader = $ws.Range("J1").Text
$titleHeader = $ws.Range("K1").Text
$deptValue = $ws.Range("J3").Text
$titleValue = $ws.Range("K3").Text

$ws.Range("H1").Value = $deptHeader
$ws.Range("I1").Value = $titleHeader
$ws.Range("H3").Value = $deptValue
$ws.Range("I3").Value = $titleValue

# --- Clear the now-vacated "Campus"/"Room" data (row5) and the old J/K columns ---
$ws.Range("H5:I5").ClearContents()
$ws.Range("J1:K1").ClearContents()
$ws.Range("J3:K3").ClearContents()

# --- Give column J (now otherwise empty) an explicit width, matching the edited layout ---
$ws.Columns("J").ColumnWidth = 9.6

# --- Update selection to reflect the edited state ---
$ws.Range("H3").Select() | Out-Null
